# Update Leve profit calculations across Sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 48534.81
$ws.Range("I129").Value = 514.5
$ws.Range("J129").Value = 78085.766
$ws.Range("K129").Value = 1543.5
$ws.Range("L129").Value = 234257.298
$ws.Range("M129").Value = 3456.5
$ws.Range("N129").Value = -244257.298

$ws.Range("H137").Value = 37038670
$ws.Range("I137").Value = 1260.8636
$ws.Range("K137").Value = 3782.5908
$ws.Range("M137").Value = -1232.5908

$ws.Range("H138").Value = 3392.0483
$ws.Range("I138").Value = 3087.3125
$ws.Range("J138").Value = 3498.0435
$ws.Range("K138").Value = 9261.9375
$ws.Range("L138").Value = 10494.1305
$ws.Range("M138").Value = -4121.9375
$ws.Range("N138").Value = -20774.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 55556612
$ws.Range("J2").Value = 1390.6364
$ws.Range("L2").Value = 1390.6364
$ws.Range("N2").Value = -1616.6364

$ws.Range("H3").Value = 8000
$ws.Range("J3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("N3").Value = -8230

$ws.Range("H32").Value = 6581.33
$ws.Range("I32").Value = 4937.033
$ws.Range("J32").Value = 21380
$ws.Range("K32").Value = 4937.033
$ws.Range("L32").Value = 21380
$ws.Range("M32").Value = -4650.033
$ws.Range("N32").Value = -21954

$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -4744

$ws.Range("H116").Value = 55556612
$ws.Range("J116").Value = 1390.6364
$ws.Range("L116").Value = 1390.6364
$ws.Range("N116").Value = -5978.6364

$ws.Range("H132").Value = 5436.83
$ws.Range("I132").Value = 5847.533
$ws.Range("J132").Value = 3126.625
$ws.Range("K132").Value = 17542.599
$ws.Range("L132").Value = 9379.875
$ws.Range("M132").Value = -15012.599
$ws.Range("N132").Value = -14439.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 55556612
$ws.Range("J3").Value = 1390.6364
$ws.Range("L3").Value = 1390.6364
$ws.Range("N3").Value = -1618.6364

$ws.Range("H99").Value = 772.1111
$ws.Range("J99").Value = 833
$ws.Range("L99").Value = 833
$ws.Range("N99").Value = -3829

$ws.Range("H126").Value = 36000
$ws.Range("J126").Value = 36000
$ws.Range("L126").Value = 36000
$ws.Range("N126").Value = -45880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3023.7742
$ws.Range("I31").Value = 2643.5715
$ws.Range("J31").Value = 3822.2
$ws.Range("K31").Value = 2643.5715
$ws.Range("L31").Value = 3822.2
$ws.Range("M31").Value = -2348.5715
$ws.Range("N31").Value = -4412.2

$ws.Range("H34").Value = 3023.7742
$ws.Range("I34").Value = 2643.5715
$ws.Range("J34").Value = 3822.2
$ws.Range("K34").Value = 2643.5715
$ws.Range("L34").Value = 3822.2
$ws.Range("M34").Value = -2441.5715
$ws.Range("N34").Value = -4226.2

$ws.Range("H94").Value = 66667628
$ws.Range("I94").Value = 142857780
$ws.Range("J94").Value = 1257
$ws.Range("K94").Value = 142857780
$ws.Range("L94").Value = 1257
$ws.Range("M94").Value = -142857329
$ws.Range("N94").Value = -2159

$ws.Range("H132").Value = 4631136.5
$ws.Range("I132").Value = 1194.6522
$ws.Range("J132").Value = 31253304
$ws.Range("K132").Value = 3583.9566
$ws.Range("L132").Value = 93759912
$ws.Range("M132").Value = -1053.9566
$ws.Range("N132").Value = -93764972

$ws.Range("H134").Value = 2036.7076
$ws.Range("I134").Value = 2036.7076
$ws.Range("K134").Value = 6110.1228
$ws.Range("M134").Value = -3575.1228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 883.6667
$ws.Range("I5").Value = 354.16666
$ws.Range("J5").Value = 3001.6667
$ws.Range("K5").Value = 1062.49998
$ws.Range("L5").Value = 9005.000100000001
$ws.Range("M5").Value = -950.4999800000001
$ws.Range("N5").Value = -9229.000100000001

$ws.Range("H122").Value = 1393.4546
$ws.Range("J122").Value = 1491
$ws.Range("L122").Value = 13419
$ws.Range("N122").Value = -18319

$ws.Range("H131").Value = 1918165.6
$ws.Range("J131").Value = 2268370.5
$ws.Range("L131").Value = 6805111.5
$ws.Range("N131").Value = -6815191.5

$ws.Range("H135").Value = 883.6667
$ws.Range("I135").Value = 354.16666
$ws.Range("J135").Value = 3001.6667
$ws.Range("K135").Value = 3187.49994
$ws.Range("L135").Value = 27015.0003
$ws.Range("M135").Value = -652.4999399999997
$ws.Range("N135").Value = -32085.0003

$ws.Range("H140").Value = 2640.8333
$ws.Range("I140").Value = 1355.7142
$ws.Range("J140").Value = 3170
$ws.Range("K140").Value = 4067.1426
$ws.Range("L140").Value = 9510
$ws.Range("M140").Value = 1112.8574
$ws.Range("N140").Value = -19870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2179.5652
$ws.Range("I97").Value = 1966.875
$ws.Range("J97").Value = 2665.7144
$ws.Range("K97").Value = 1966.875
$ws.Range("L97").Value = 2665.7144
$ws.Range("M97").Value = -1470.875
$ws.Range("N97").Value = -3657.7144

$ws.Range("H113").Value = 50001084
$ws.Range("I113").Value = 125000504
$ws.Range("J113").Value = 1466.6666
$ws.Range("K113").Value = 125000504
$ws.Range("L113").Value = 1466.6666
$ws.Range("M113").Value = -124998334
$ws.Range("N113").Value = -5806.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1433.84
$ws.Range("I46").Value = 829
$ws.Range("J46").Value = 1549.0476
$ws.Range("K46").Value = 829
$ws.Range("L46").Value = 1549.0476
$ws.Range("M46").Value = -641
$ws.Range("N46").Value = -1925.0476

$ws.Range("H122").Value = 4112.864
$ws.Range("I122").Value = 5215.25
$ws.Range("J122").Value = 2790
$ws.Range("K122").Value = 15645.75
$ws.Range("L122").Value = 8370
$ws.Range("M122").Value = -13195.75
$ws.Range("N122").Value = -13270

$ws.Range("H133").Value = 25819.385
$ws.Range("J133").Value = 25819.385
$ws.Range("L133").Value = 25819.385
$ws.Range("N133").Value = -30879.385

$ws.Range("H136").Value = 3522.0378
$ws.Range("I136").Value = 3516.6956
$ws.Range("J136").Value = 3557.1428
$ws.Range("K136").Value = 10550.0868
$ws.Range("L136").Value = 10671.4284
$ws.Range("M136").Value = -8000.086800000001
$ws.Range("N136").Value = -15771.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 267272.72
$ws.Range("J54").Value = 44000
$ws.Range("L54").Value = 44000
$ws.Range("N54").Value = -45040
